{"js": "// Prueba tecnica IAS Skills segundo cambio\n// Fill in the first empty line after the last existing query with the\n// \"CREATE or REPLACE VIEW ...\" statement, then add two more paragraphs\n// for the rest of the view definition (select ... / from canciones;).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the last paragraph of the previous query (unique anchor text) so\n// we land on the correct blank line even if blank-paragraph counts drift.\nconst anchorText =\n  \"INNER JOIN artistas on canciones.artistas_fk = artistas.id;\";\n\nlet anchorIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === anchorText) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find anchor paragraph for the new SQL view.\");\n}\n\n// The new statement goes on the 6th blank paragraph following the anchor.\nconst targetIndex = anchorIndex + 6;\nconst targetParagraph = paragraphs.items[targetIndex];\n\ntargetParagraph.insertText(\n  \"CREATE or REPLACE VIEW Canciones_y_duraci\u00f3n as \",\n  \"Replace\"\n);\n\nconst secondParagraph = targetParagraph.insertParagraph(\n  \"select canciones.titulo as titulo_cancion, canciones.duracion as duracion_cancion\",\n  \"After\"\n);\n\nsecondParagraph.insertParagraph(\"from canciones;\", \"After\");\n\nawait context.sync();\n", "ps1": "# Prueba tecnica IAS Skills segundo cambio\n# Fill in the first empty line after the last existing query with the\n# \"CREATE or REPLACE VIEW ...\" statement, then add two more paragraphs\n# for the rest of the view definition (select ... / from canciones;).\n\n$d = $word.ActiveDocument\n$paragraphs = $d.Paragraphs\n\n# Locate the last paragraph of the previous query (unique anchor text) so\n# we land on the correct blank line even if blank-paragraph counts drift.\n$anchorText = \"INNER JOIN artistas on canciones.artistas_fk = artistas.id;\"\n$anchorIndex = -1\nfor ($i = 1; $i -le $paragraphs.Count; $i++) {\n    $text = $paragraphs.Item($i).Range.Text.TrimEnd()\n    if ($text -eq $anchorText) {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not find anchor paragraph for the new SQL view.\"\n}\n\n# The new statement goes on the 6th blank paragraph following the anchor.\n$targetIndex = $anchorIndex + 6\n$target = $d.Paragraphs.Item($targetIndex)\n$target.Range.Text = \"CREATE or REPLACE VIEW Canciones_y_duraci\u00f3n as \"\n$target.Range.InsertParagraphAfter()\n\n$second = $d.Paragraphs.Item($targetIndex + 1)\n$second.Range.Text = \"select canciones.titulo as titulo_cancion, canciones.duracion as duracion_cancion\"\n$second.Range.InsertParagraphAfter()\n\n$third = $d.Paragraphs.Item($targetIndex + 2)\n$third.Range.Text = \"from canciones;\"\n"}
